# Apply scheduled-runner market-data refresh to the Leve profit tables.
# For each changed leve row, currentAveragePrice* (H/I/J/K), LevePrice* (L),
# and LeveProfit* (M/N) columns are rewritten with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 277.02856
$ws.Range("I6").Value = 129.5
$ws.Range("K6").Value = 388.5
$ws.Range("M6").Value = -276.5

$ws.Range("H32").Value = 1286.2858
$ws.Range("J32").Value = 1472.9
$ws.Range("L32").Value = 1472.9
$ws.Range("N32").Value = -2124.9

$ws.Range("H48").Value = 4998.5
$ws.Range("J48").Value = 4998.5
$ws.Range("L48").Value = 14995.5
$ws.Range("N48").Value = -15579.5

$ws.Range("H56").Value = 4998.5
$ws.Range("J56").Value = 4998.5
$ws.Range("L56").Value = 14995.5
$ws.Range("N56").Value = -16063.5

$ws.Range("H120").Value = 125000
$ws.Range("J120").Value = 125000
$ws.Range("L120").Value = 125000
$ws.Range("N120").Value = -134676

$ws.Range("H129").Value = 66670920
$ws.Range("J129").Value = 14000
$ws.Range("L129").Value = 42000
$ws.Range("N129").Value = -52000

$ws.Range("H138").Value = 2931.7676
$ws.Range("J138").Value = 3630.7183
$ws.Range("L138").Value = 10892.1549
$ws.Range("N138").Value = -21172.1549

$ws.Range("H141").Value = 2585.9678
$ws.Range("I141").Value = 1008.04
$ws.Range("J141").Value = 9160.666999999999
$ws.Range("K141").Value = 3024.12
$ws.Range("L141").Value = 27482.001
$ws.Range("M141").Value = 2155.88
$ws.Range("N141").Value = -37842.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3970.7
$ws.Range("I32").Value = 3075.3618
$ws.Range("K32").Value = 3075.3618
$ws.Range("M32").Value = -2788.3618

$ws.Range("H61").Value = 9236
$ws.Range("I61").Value = 9281.471
$ws.Range("K61").Value = 9281.471
$ws.Range("M61").Value = -9069.471

$ws.Range("H113").Value = 73000
$ws.Range("J113").Value = 73000
$ws.Range("L113").Value = 73000
$ws.Range("N113").Value = -81678

$ws.Range("H132").Value = 2487.94
$ws.Range("I132").Value = 2275.889
$ws.Range("K132").Value = 6827.667
$ws.Range("M132").Value = -4297.667

$ws.Range("H135").Value = 73364.875
$ws.Range("J135").Value = 73364.875
$ws.Range("L135").Value = 73364.875
$ws.Range("N135").Value = -83504.875

$ws.Range("H136").Value = 9236
$ws.Range("I136").Value = 9281.471
$ws.Range("K136").Value = 27844.413
$ws.Range("M136").Value = -25294.413

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4987.5654
$ws.Range("I105").Value = 4032.125
$ws.Range("J105").Value = 7171.4287
$ws.Range("K105").Value = 4032.125
$ws.Range("L105").Value = 7171.4287
$ws.Range("M105").Value = -2285.125
$ws.Range("N105").Value = -10665.4287

$ws.Range("H134").Value = 4051.7646
$ws.Range("I134").Value = 4025.8572
$ws.Range("J134").Value = 4172.6665
$ws.Range("K134").Value = 12077.5716
$ws.Range("L134").Value = 12517.9995
$ws.Range("M134").Value = -9542.571599999999
$ws.Range("N134").Value = -17587.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 482.58334
$ws.Range("I22").Value = 301.88235
$ws.Range("K22").Value = 301.88235
$ws.Range("M22").Value = 48.11765000000003

$ws.Range("H31").Value = 2954.9355
$ws.Range("J31").Value = 3786.2666
$ws.Range("L31").Value = 3786.2666
$ws.Range("N31").Value = -4376.2666

$ws.Range("H34").Value = 2954.9355
$ws.Range("J34").Value = 3786.2666
$ws.Range("L34").Value = 3786.2666
$ws.Range("N34").Value = -4190.2666

$ws.Range("H58").Value = 1665.6154
$ws.Range("I58").Value = 1415.9
$ws.Range("K58").Value = 1415.9
$ws.Range("M58").Value = -1212.9

$ws.Range("H86").Value = 7002.25
$ws.Range("I86").Value = 6092.7
$ws.Range("K86").Value = 6092.7
$ws.Range("M86").Value = -4969.7

$ws.Range("H89").Value = 7002.25
$ws.Range("I89").Value = 6092.7
$ws.Range("K89").Value = 30463.5
$ws.Range("M89").Value = -24847.5

$ws.Range("H99").Value = 6202.25
$ws.Range("I99").Value = 6055.4736
$ws.Range("J99").Value = 6760
$ws.Range("K99").Value = 6055.4736
$ws.Range("L99").Value = 6760
$ws.Range("M99").Value = -4557.4736
$ws.Range("N99").Value = -9756

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H125").Value = 81990
$ws.Range("J125").Value = 81990
$ws.Range("L125").Value = 81990
$ws.Range("N125").Value = -86910

$ws.Range("H126").Value = 6202.25
$ws.Range("I126").Value = 6055.4736
$ws.Range("J126").Value = 6760
$ws.Range("K126").Value = 18166.4208
$ws.Range("L126").Value = 20280
$ws.Range("M126").Value = -15696.4208
$ws.Range("N126").Value = -25220

$ws.Range("H134").Value = 3740.4075
$ws.Range("I134").Value = 1174.4375
$ws.Range("J134").Value = 7472.727
$ws.Range("K134").Value = 3523.3125
$ws.Range("L134").Value = 22418.181
$ws.Range("M134").Value = -988.3125
$ws.Range("N134").Value = -27488.181

$ws.Range("H136").Value = 1665.6154
$ws.Range("I136").Value = 1415.9
$ws.Range("K136").Value = 4247.700000000001
$ws.Range("M136").Value = -1697.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2781.75
$ws.Range("I14").Value = 2781.75
$ws.Range("K14").Value = 8345.25
$ws.Range("M14").Value = -8172.25

$ws.Range("H38").Value = 14490.143
$ws.Range("J38").Value = 33116.332
$ws.Range("L38").Value = 99348.99600000001
$ws.Range("N38").Value = -100042.996

$ws.Range("H80").Value = 6088.5
$ws.Range("I80").Value = 4477
$ws.Range("K80").Value = 13431
$ws.Range("M80").Value = -12495

$ws.Range("H83").Value = 6088.5
$ws.Range("I83").Value = 4477
$ws.Range("K83").Value = 40293
$ws.Range("M83").Value = -35613

$ws.Range("H132").Value = 3737.0386
$ws.Range("J132").Value = 3976.913
$ws.Range("L132").Value = 35792.217
$ws.Range("N132").Value = -40852.217

$ws.Range("H137").Value = 4186.9473
$ws.Range("J137").Value = 3386
$ws.Range("L137").Value = 10158
$ws.Range("N137").Value = -20358

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 46333.332
$ws.Range("J93").Value = 46333.332
$ws.Range("L93").Value = 46333.332
$ws.Range("N93").Value = -50077.332

$ws.Range("H122").Value = 2777.25
$ws.Range("I122").Value = 2565
$ws.Range("K122").Value = 7695
$ws.Range("M122").Value = -5245

$ws.Range("H128").Value = 100000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 100000
$ws.Range("K128").Value = 0
$ws.Range("M128").Value = 100000
$ws.Range("N128").Value = -109960

$ws.Range("H132").Value = 8201.143
$ws.Range("I132").Value = 8900.166999999999
$ws.Range("K132").Value = 26700.501
$ws.Range("M132").Value = -24170.501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3597.6924
$ws.Range("I132").Value = 3508.2273
$ws.Range("J132").Value = 3713.4707
$ws.Range("K132").Value = 10524.6819
$ws.Range("L132").Value = 11140.4121
$ws.Range("M132").Value = -7994.6819
$ws.Range("N132").Value = -16200.4121

$ws.Range("H136").Value = 4718.1465
$ws.Range("I136").Value = 4307.2354
$ws.Range("K136").Value = 12921.7062
$ws.Range("M136").Value = -10371.7062

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 28611
$ws.Range("J47").Value = 28611
$ws.Range("L47").Value = 28611
$ws.Range("N47").Value = -29755

$ws.Range("H100").Value = 2762.1667
$ws.Range("I100").Value = 2701.25
$ws.Range("J100").Value = 3249.5
$ws.Range("K100").Value = 5402.5
$ws.Range("L100").Value = 6499
$ws.Range("M100").Value = -4861.5
$ws.Range("N100").Value = -7581

$ws.Range("H126").Value = 2768.125
$ws.Range("I126").Value = 2418
$ws.Range("K126").Value = 7254
$ws.Range("M126").Value = -4784

$ws.Range("H132").Value = 1693.7059
$ws.Range("I132").Value = 1595.8889
$ws.Range("K132").Value = 4787.6667
$ws.Range("M132").Value = -2257.6667

$ws.Range("H136").Value = 11591.263
$ws.Range("I136").Value = 11957.5
$ws.Range("K136").Value = 35872.5
$ws.Range("M136").Value = -33322.5
